$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "59.121.73"
Set-TextValue "E2" "  +3.26%  "
Set-TextValue "D3" "2.592.58"
Set-TextValue "E3" "  +1.43%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "521.76"
Set-TextValue "E5" "  +1.06%  "
Set-TextValue "D6" "139.71"
Set-TextValue "E6" "  -1.29%  "
Set-TextValue "E7" "  -0.10%  "
Set-TextValue "D8" "0.565"
Set-TextValue "E8" "  +0.78%  "
Set-TextValue "D9" "2.606.15"
Set-TextValue "E9" "  +1.52%  "
Set-TextValue "D10" "6.54"
Set-TextValue "E10" "  -0.47%  "
Set-TextValue "D11" "0.101"
Set-TextValue "E11" "  +1.40%  "
Set-TextValue "D12" "0.330"
Set-TextValue "E12" "  +2.19%  "
Set-TextValue "E13" "  +3.26%  "
Set-TextValue "D14" "3.051.81"
Set-TextValue "E14" "  +1.42%  "
Set-TextValue "D15" "59.035.03"
Set-TextValue "E15" "  +3.14%  "
Set-TextValue "D16" "20.45"
Set-TextValue "E16" "  +1.75%  "
Set-TextValue "D17" "2.595.37"
Set-TextValue "E17" "  +1.59%  "
Set-TextValue "D18" "0.0000133"
Set-TextValue "E18" "  +0.50%  "
Set-TextValue "D19" "339.39"
Set-TextValue "E19" "  +1.55%  "
Set-TextValue "D20" "4.32"
Set-TextValue "E20" "  +1.33%  "
Set-TextValue "D21" "10.14"
Set-TextValue "E21" "  +0.17%  "
Set-TextValue "D22" "6.46"
Set-TextValue "E22" "  +3.40%  "
Set-TextValue "D23" "1.00"
Set-TextValue "E23" "  +0.02%  "
Set-TextValue "D24" "66.40"
Set-TextValue "E24" "  +1.85%  "
Set-TextValue "D25" "0.168"
Set-TextValue "E25" "  +2.10%  "
Set-TextValue "D26" "0.404"
Set-TextValue "E26" "  +1.22%  "
Set-TextValue "E27" "  +0.02%  "
Set-TextValue "D28" "7.02"
Set-TextValue "E28" "  +1.52%  "
Set-TextValue "D29" "0.998"
Set-TextValue "E29" "  +0.05%  "
Set-TextValue "D30" "0.0₃0728"
Set-TextValue "E30" "  -1.59%  "
Set-TextValue "D31" "5.95"
Set-TextValue "E31" "  -6.70%  "
Set-TextValue "D32" "1.58"
Set-TextValue "E32" "  +0.14%  "
Set-TextValue "D33" "18.79"
Set-TextValue "E33" "  +1.35%  "
Set-TextValue "D34" "148.86"
Set-TextValue "E34" "  -0.36%  "
Set-TextValue "D35" "4.00"
Set-TextValue "E35" "  +0.74%  "
Set-TextValue "E36" "  -0.81%  "
Set-TextValue "D37" "36.75"
Set-TextValue "E37" "  +2.47%  "
Set-TextValue "D38" "1.47"
Set-TextValue "E38" "  +2.52%  "
Set-TextValue "D39" "0.831"
Set-TextValue "E39" "  +0.69%  "
Set-TextValue "D40" "0.811"
Set-TextValue "E40" "  -6.31%  "
Set-TextValue "D41" "3.52"
Set-TextValue "E41" "  +0.25%  "
Set-TextValue "E42" "  -0.05%  "
Set-TextValue "D43" "274.38"
Set-TextValue "E43" "  +1.89%  "
Set-TextValue "E44" "  +0.93%  "
Set-TextValue "D45" "0.592"
Set-TextValue "E45" "  +1.48%  "
Set-TextValue "D46" "0.0951"
Set-TextValue "E46" "  +0.18%  "
Set-TextValue "D47" "0.0519"
Set-TextValue "E47" "  +0.03%  "
Set-TextValue "D48" "18.56"
Set-TextValue "E48" "  -0.57%  "
Set-TextValue "D49" "1.971.34"
Set-TextValue "E49" "  +0.50%  "
Set-TextValue "D50" "4.55"
Set-TextValue "E50" "  +0.88%  "
Set-TextValue "D51" "0.0219"
Set-TextValue "E51" "  +0.78%  "
